$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - qty upto date
$ws.Range("C8").Value = 90

# Row 9 - qty upto date + upto date amount (94 * 256 = 24064.00)
$ws.Range("C9").Value = 94
$ws.Range("G9").Value = "'24064.00"

# Row 10 - qty upto date + upto date amount (54 * 472 = 25488.00)
$ws.Range("C10").Value = 54
$ws.Range("G10").Value = "'25488.00"

# Row 11 - qty upto date + upto date amount (64 * 662 = 42368.00)
$ws.Range("C11").Value = 64
$ws.Range("G11").Value = "'42368.00"

# Row 12 - qty upto date
$ws.Range("C12").Value = 90

# Row 13 - qty upto date + upto date amount (11 * 136 = 1496.00)
$ws.Range("C13").Value = 11
$ws.Range("G13").Value = "'1496.00"

# Row 14 - qty upto date + upto date amount (20 * 23 = 460.00)
$ws.Range("C14").Value = 20
$ws.Range("G14").Value = "'460.00"

# Row 15 - qty upto date
$ws.Range("C15").Value = 33

# Row 16 - qty upto date
$ws.Range("C16").Value = 50

# Row 17 - qty upto date
$ws.Range("C17").Value = 50

# Row 19 - Grand Total Rs.
$ws.Range("G19").Value = "'93876.00"
$ws.Range("H19").Value = "'93876.00"

# Row 21 - NET PAYABLE AMOUNT Rs.
$ws.Range("G21").Value = "'93876.00"
$ws.Range("H21").Value = "'93876.00"
